$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1 / sheetId 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 343
$ws1.Range("F5").Value = 5073
$ws1.Range("F9").Value = 765
$ws1.Range("F10").Value = 249
$ws1.Range("F11").Value = 7

# Sheet "全部类型" (index 4 / sheetId 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 343
$ws4.Range("F5").Value = 5073
$ws4.Range("F9").Value = 765
$ws4.Range("F11").Value = 249
$ws4.Range("F12").Value = 7
